# Add two new "Battling" sub-bullets ("Abilities" and "Items") right
# after the existing "Non-attack damage" bullet (and before the next
# bullet, a lone "-") inside the System/Battling outline.
#
# Non-attack damage
# Abilities      <- new
# Items          <- new
# -

$d = $word.ActiveDocument

# Locate the "Non-attack damage" paragraph by its text so the edit is
# resilient to any paragraph-numbering shifts.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Non-attack damage") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find the 'Non-attack damage' paragraph"
}

# Insert "Items" first, then "Abilities" - both are inserted
# immediately after the anchor, so inserting in reverse order leaves
# them in the desired final order (Abilities, then Items).
$r = $anchor.Range
$r.InsertParagraphAfter()
$itemsPara = $anchor.Next()
$itemsPara.Range.Text = "Items"

$r = $anchor.Range
$r.InsertParagraphAfter()
$abilitiesPara = $anchor.Next()
$abilitiesPara.Range.Text = "Abilities"
